# Generate Report for Handoff
# Refresh the handoff/report timestamps and priority flag for the
# 02ed08f4-b7c3-4428-aec1-eb3ef9c4bb42.md row (rows 7-12 across the three
# sheets, which all mirror the same source row) to reflect the newly
# generated handoff.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("G7:G12").Value = "2016-09-06 04:26:56"

# zh-cn sheet: "Priority" column (E) and "Latest Handoff Datetime" column (H)
$wsZhCn.Range("E7:E12").Value = "ht"
$wsZhCn.Range("H7:H12").Value = "2016-09-06 04:26:47"

# de-de sheet: "Priority" column (E) and "Latest Handoff Datetime" column (H)
$wsDeDe.Range("E7:E12").Value = "ht"
$wsDeDe.Range("H7:H12").Value = "2016-09-06 04:26:56"
